# Restructure the "Correspondence address" block:
#   - Para 1 ("Correspondence address:") loses the _GoBack bookmark.
#   - Para 2 ("{ownerName}") moves from ind=-58 to ind=893 and drops the
#     19 leading spaces in front of "{".
#   - Para 3 ("{@ownerAddress}") moves from ind=-58 to ind=893 and gains
#     the _GoBack bookmark at its start.

$d = $word.ActiveDocument

$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the three paragraphs by their text content (robust to index drift).
$paraAddress = $null
$paraOwnerName = $null
$paraOwnerAddress = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Correspondence address:*") {
        $paraAddress = $p
    } elseif ($t -like "*{ownerName}*") {
        $paraOwnerName = $p
    } elseif ($t -like "*{@ownerAddress}*") {
        $paraOwnerAddress = $p
    }
}

# --- Paragraph 1: "Correspondence address:" — strip the _GoBack bookmark ---
$xml1 = "<w:p $W><w:pPr><w:spacing w:after=`"430`" w:line=`"240`" w:lineRule=`"auto`"/><w:ind w:left=`"-58`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr><w:t>Correspondence address:</w:t></w:r></w:p>"
$paraAddress.Range.InsertXML($xml1)

# --- Paragraph 2: "{ownerName}" — re-indent, drop the leading spaces ---
$xml2 = "<w:p $W><w:pPr><w:spacing w:after=`"0`"/><w:ind w:left=`"893`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"FFFFFF`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr><w:t>ownerName</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr><w:t>}</w:t></w:r></w:p>"
$paraOwnerName.Range.InsertXML($xml2)

# --- Paragraph 3: "{@ownerAddress}" — re-indent, gain the _GoBack bookmark ---
$xml3 = "<w:p $W><w:pPr><w:spacing w:after=`"0`" w:line=`"360`" w:lineRule=`"auto`"/><w:ind w:left=`"893`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr></w:pPr><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr><w:t>{@</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr><w:t>ownerAddress</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b/><w:sz w:val=`"20`"/><w:szCs w:val=`"20`"/></w:rPr><w:t>}</w:t></w:r></w:p>"
$paraOwnerAddress.Range.InsertXML($xml3)

Write-Output "done"
